$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 13 so the existing "frazier"/"yngve" rows
# (currently rows 14 and 15) shift down to rows 16 and 17, leaving row 13
# free for the new "schcount" row (and row 14 stays empty, same as the
# blank row 13 that existed before the edit).
$ws.Range("13:14").Insert() | Out-Null
$ws.Range("13:14").Clear() | Out-Null

# New row 13: schcount
$ws.Range("A13").Value = "schcount"

$ws.Range("B13").Value = 1
$ws.Range("B13").Style = "Good"
$ws.Range("C13").Value = [double]"5.8983000000000003E-7"
$ws.Range("C13").NumberFormat = "0.00E+00"

$ws.Range("D13").Value = 1
$ws.Range("D13").Style = "Good"
$ws.Range("E13").Value = [double]"1.6000000000000001E-3"

$ws.Range("F13").Value = 1
$ws.Range("F13").Style = "Good"
$ws.Range("G13").Value = [double]"5.3E-3"

$ws.Range("H13").Value = 1
$ws.Range("H13").Style = "Good"
$ws.Range("I13").Value = [double]"1.5543999999999999E-4"
$ws.Range("I13").NumberFormat = "0.00E+00"

$ws.Range("J13").Value = 1
$ws.Range("J13").Style = "Good"
$ws.Range("K13").Value = [double]"6.016E-5"
$ws.Range("K13").NumberFormat = "0.00E+00"

$ws.Range("L6").Select() | Out-Null
